$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 6010
$ws.Range("I29").Value = 276
$ws.Range("J29").Value = 9832.666999999999
$ws.Range("K29").Value = 828
$ws.Range("L29").Value = 29498.001
$ws.Range("M29").Value = -547
$ws.Range("N29").Value = -30060.001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 9700
$ws.Range("J46").Value = 7400
$ws.Range("L46").Value = 22200
$ws.Range("N46").Value = -22438

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 577.3570999999999
$ws.Range("I53").Value = 573.5833
$ws.Range("J53").Value = 600
$ws.Range("K53").Value = 573.5833
$ws.Range("L53").Value = 600
$ws.Range("M53").Value = 63.41669999999999
$ws.Range("N53").Value = -1874

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H60").Value = 9700
$ws.Range("J60").Value = 7400
$ws.Range("L60").Value = 22200
$ws.Range("N60").Value = -23168

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1572.4286
$ws.Range("I86").Value = 1600.6
$ws.Range("J86").Value = 1502
$ws.Range("K86").Value = 1600.6
$ws.Range("L86").Value = 1502
$ws.Range("M86").Value = -477.5999999999999
$ws.Range("N86").Value = -3748

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1572.4286
$ws.Range("I89").Value = 1600.6
$ws.Range("J89").Value = 1502
$ws.Range("K89").Value = 8003
$ws.Range("L89").Value = 7510
$ws.Range("M89").Value = -2387
$ws.Range("N89").Value = -18742

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1948.6
$ws.Range("I98").Value = 2185.75
$ws.Range("K98").Value = 2185.75
$ws.Range("M98").Value = -687.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1948.6
$ws.Range("I122").Value = 2185.75
$ws.Range("K122").Value = 6557.25
$ws.Range("M122").Value = -4107.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1500
$ws.Range("I137").Value = 1500
$ws.Range("K137").Value = 4500
$ws.Range("M137").Value = -1950

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4259.9
$ws.Range("I32").Value = 3955.4443
$ws.Range("K32").Value = 3955.4443
$ws.Range("M32").Value = -3668.4443

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 800
$ws.Range("I74").Value = 800
$ws.Range("K74").Value = 800
$ws.Range("M74").Value = 74

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 52000
$ws.Range("J76").Value = 52000
$ws.Range("L76").Value = 52000
$ws.Range("N76").Value = -52676

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 800
$ws.Range("I77").Value = 800
$ws.Range("K77").Value = 4000
$ws.Range("M77").Value = 368

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 52000
$ws.Range("J79").Value = 52000
$ws.Range("L79").Value = 52000
$ws.Range("N79").Value = -54340

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 466.66666
$ws.Range("I132").Value = 466.66666
$ws.Range("K132").Value = 1399.99998
$ws.Range("M132").Value = 1130.00002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2741.2
$ws.Range("I86").Value = 2741.2
$ws.Range("K86").Value = 2741.2
$ws.Range("M86").Value = -1618.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2741.2
$ws.Range("I89").Value = 2741.2
$ws.Range("K89").Value = 13706
$ws.Range("M89").Value = -8090

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2750
$ws.Range("I134").Value = 2750
$ws.Range("K134").Value = 8250
$ws.Range("M134").Value = -5715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1666.6666
$ws.Range("I31").Value = 1666.6666
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1666.6666
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1371.6666
$ws.Range("N31").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1666.6666
$ws.Range("I34").Value = 1666.6666
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1666.6666
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1464.6666
$ws.Range("N34").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 17049.834
$ws.Range("J92").Value = 17049.834
$ws.Range("L92").Value = 17049.834
$ws.Range("N92").Value = -22041.834

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5378
$ws.Range("I132").Value = 4504
$ws.Range("K132").Value = 13512
$ws.Range("M132").Value = -10982

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1654.4445
$ws.Range("I121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("M121").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2100.6
$ws.Range("I102").Value = 2000.875
$ws.Range("J102").Value = 2499.5
$ws.Range("K102").Value = 2000.875
$ws.Range("L102").Value = 2499.5
$ws.Range("M102").Value = -378.875
$ws.Range("N102").Value = -5743.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2333.3333
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6024.875
$ws.Range("I40").Value = 4186.143
$ws.Range("K40").Value = 4186.143
$ws.Range("M40").Value = -4050.143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 11285.857
$ws.Range("J46").Value = 17500
$ws.Range("L46").Value = 17500
$ws.Range("N46").Value = -17876

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3000
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3000
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1450
$ws.Range("I82").Value = 1316.6666
$ws.Range("J82").Value = 1650
$ws.Range("K82").Value = 1316.6666
$ws.Range("L82").Value = 1650
$ws.Range("M82").Value = -955.6666
$ws.Range("N82").Value = -2372

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1450
$ws.Range("I85").Value = 1316.6666
$ws.Range("J85").Value = 1650
$ws.Range("K85").Value = 1316.6666
$ws.Range("L85").Value = 1650
$ws.Range("M85").Value = -68.66660000000002
$ws.Range("N85").Value = -4146

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 999.5
$ws.Range("I93").Value = 999.5
$ws.Range("K93").Value = 999.5
$ws.Range("M93").Value = 248.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 27100
$ws.Range("I136").Value = 28500
$ws.Range("J136").Value = 25000
$ws.Range("K136").Value = 85500
$ws.Range("L136").Value = 75000
$ws.Range("M136").Value = -82950
$ws.Range("N136").Value = -80100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 74250
$ws.Range("J103").Value = 74250
$ws.Range("L103").Value = 74250
$ws.Range("N103").Value = -76594

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1613.1666
$ws.Range("I107").Value = 1260
$ws.Range("J107").Value = 1966.3334
$ws.Range("K107").Value = 3780
$ws.Range("L107").Value = 5899.0002
$ws.Range("M107").Value = -1860
$ws.Range("N107").Value = -9739.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1420
$ws.Range("I132").Value = 1456.6666
$ws.Range("K132").Value = 4369.9998
$ws.Range("M132").Value = -1839.9998
